$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 317 (pushes existing rows 317-338 down to 319-340)
$ws.Rows.Item(317).Insert()
$ws.Rows.Item(317).Insert()

# New row 317: Apio, Primera, week of 44615
$ws.Cells.Item(317, 1).Value = 8
$ws.Cells.Item(317, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(317, 3).Value = "Coquimbo"
$ws.Cells.Item(317, 4).Value = 44615
$ws.Cells.Item(317, 5).Value = 4
$ws.Cells.Item(317, 6).Value = 100112017
$ws.Cells.Item(317, 7).Value = "Apio"
$ws.Cells.Item(317, 8).Value = "Americana (o)"
$ws.Cells.Item(317, 9).Value = "Primera"
$ws.Cells.Item(317, 10).Value = 2200
$ws.Cells.Item(317, 11).Value = 8000
$ws.Cells.Item(317, 12).Value = 9000
$ws.Cells.Item(317, 13).Value = 8500
$ws.Cells.Item(317, 14).Value = "$/docena de matas"
$ws.Cells.Item(317, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(317, 16).Value = 1417
$ws.Cells.Item(317, 17).Value = 6
$ws.Cells.Item(317, 18).Value = "Hortaliza"

# New row 318: Apio, Segunda, week of 44615
$ws.Cells.Item(318, 1).Value = 8
$ws.Cells.Item(318, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(318, 3).Value = "Coquimbo"
$ws.Cells.Item(318, 4).Value = 44615
$ws.Cells.Item(318, 5).Value = 4
$ws.Cells.Item(318, 6).Value = 100112017
$ws.Cells.Item(318, 7).Value = "Apio"
$ws.Cells.Item(318, 8).Value = "Americana (o)"
$ws.Cells.Item(318, 9).Value = "Segunda"
$ws.Cells.Item(318, 10).Value = 1300
$ws.Cells.Item(318, 11).Value = 6000
$ws.Cells.Item(318, 12).Value = 7000
$ws.Cells.Item(318, 13).Value = 6500
$ws.Cells.Item(318, 14).Value = "$/docena de matas"
$ws.Cells.Item(318, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(318, 16).Value = 1083
$ws.Cells.Item(318, 17).Value = 6
$ws.Cells.Item(318, 18).Value = "Hortaliza"
